# Generate Report for Archive
# Update localization status from "Ready for handoff" to "In Translation"
# across the Overview sheet (zh-cn / de-de columns) and the per-locale
# status sheets, then refresh the affected column widths to reflect the
# shorter status text (as Excel does when the report is regenerated).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: zh-cn (col E) / de-de (col F) status columns ---
$lastRow = $overview.Cells.Item($overview.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    if ($overview.Cells.Item($r, 5).Value() -eq $oldStatus) {
        $overview.Cells.Item($r, 5).Value = $newStatus
    }
    if ($overview.Cells.Item($r, 6).Value() -eq $oldStatus) {
        $overview.Cells.Item($r, 6).Value = $newStatus
    }
}

# --- Per-locale sheets: Status column (col C) ---
foreach ($ws in @($zhcn, $dede)) {
    $lastRowLoc = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    for ($r = 2; $r -le $lastRowLoc; $r++) {
        if ($ws.Cells.Item($r, 3).Value() -eq $oldStatus) {
            $ws.Cells.Item($r, 3).Value = $newStatus
        }
    }
}

# --- Resize the status columns to fit the new, shorter text ---
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
